# Regenerate the "K" column (column G) values in the save_data sheet.
# These values are recomputed upstream (K replacing the old Strike# metric,
# with std/mean regenerated and s_vals recalculated) and simply need to be
# written back into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 1
    8  = 3
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 1
    15 = 1
    16 = 1
    17 = 3
    18 = 0
    19 = 1
    20 = 2
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 3
    26 = 2
    27 = 0
    28 = 4
    29 = 2
    30 = 3
    31 = 2
    32 = 2
    33 = 2
    34 = 2
    35 = 0
    36 = 1
    37 = 4
    38 = 2
    39 = 2
    40 = 4
    41 = 2
    42 = 0
    43 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
